$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "24/10/2025"
$ws.Range("B10").Value = "Leeds"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = "West Ham"
$ws.Range("F10").Value = "L"
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0.63
$ws.Range("L10").Value = 1.53
$ws.Range("M10").Value = 9
$ws.Range("N10").Value = 13
$ws.Range("O10").Value = 3
$ws.Range("P10").Value = 5
